$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row -> remove it entirely, shifting everything below up by one
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# K2 / L2 short description for the root Extension element
$elements.Range("K2").Value = "Network Id"
$elements.Range("L2").Value = "Customer-specific identifier of the patient provider network in which the member is enrolled"
